$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.14000000000049
$ws.Range("H2").Value = 0.00002588805161851671
$ws.Range("I2").Value = 0.00002588805161851671
$ws.Range("L2").Value = 43.60444384246517
$ws.Range("M2").Value = "[21.364429568137965, 65.84445811679238]"
$ws.Range("N2").Value = 0.0002733537664807439
$ws.Range("O2").Value = 0.0002733537664807439
$ws.Range("P2").Value = 1.201289683354656
$ws.Range("Q2").Value = "[0.6352369529781168, 1.7673424137311944]"
$ws.Range("R2").Value = 0.00009828815897416909
$ws.Range("S2").Value = 0.00009828815897416909
$ws.Range("T2").Value = 62.01283808085582
$ws.Range("U2").Value = "[49.927848215860266, 74.09782794585138]"
$ws.Range("V2").Value = 0.0000000000001834088436680759
$ws.Range("W2").Value = 0.0000000000001834088436680759
$ws.Range("X2").Value = 20.33345345345385
$ws.Range("Y2").Value = 18.06858858858894
$ws.Range("Z2").Value = 22.59831831831876

# Row 3
$ws.Range("F3").Value = 25.14000000000049
$ws.Range("H3").Value = 0.00002610500728206233
$ws.Range("I3").Value = 0.00002610500728206233
$ws.Range("L3").Value = 44.5017367341492
$ws.Range("M3").Value = "[24.67775026666503, 64.32572320163337]"
$ws.Range("N3").Value = 0.00004441794900400531
$ws.Range("O3").Value = 0.00004441794900400531
$ws.Range("P3").Value = 1.465447624197041
$ws.Range("Q3").Value = "[0.9245527929483472, 2.0063424554457345]"
$ws.Range("R3").Value = 0.000001983180353271052
$ws.Range("S3").Value = 0.000001983180353271052
$ws.Range("T3").Value = 49.1522728155575
$ws.Range("U3").Value = "[36.691291239273546, 61.61325439184146]"
$ws.Range("V3").Value = 0.0000000004150482180165227
$ws.Range("W3").Value = 0.0000000004150482180165227
$ws.Range("X3").Value = 19.27651651651689
$ws.Range("Y3").Value = 17.11231231231264
$ws.Range("Z3").Value = 21.44072072072114

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 25.14000000000049
$ws.Range("H4").Value = 0.00000005000563685353399
$ws.Range("I4").Value = 0.00000005000563685353399
$ws.Range("L4").Value = 53.64957137741292
$ws.Range("M4").Value = "[33.92779486894331, 73.37134788588253]"
$ws.Range("N4").Value = 0.000001839657850499776
$ws.Range("O4").Value = 0.000001839657850499776
$ws.Range("P4").Value = 0.5220264069028087
$ws.Range("Q4").Value = "[0.14465791998511524, 0.8993948938205021]"
$ws.Range("R4").Value = 0.007780180144131288
$ws.Range("S4").Value = 0.007780180144131288
$ws.Range("T4").Value = 61.24964402408007
$ws.Range("U4").Value = "[50.65746444169683, 71.84182360646331]"
$ws.Range("V4").Value = 0.000000000000003552713678800501
$ws.Range("W4").Value = 0.000000000000003552713678800501
$ws.Range("X4").Value = 23.05129129129174
$ws.Range("Y4").Value = 21.5413813813818
$ws.Range("Z4").Value = 24.56120120120168

# Row 5
$ws.Range("F5").Value = 25.14000000000049
$ws.Range("H5").Value = 0.01025216557263364
$ws.Range("I5").Value = 0.01025216557263364
$ws.Range("L5").Value = 29.39829065925683
$ws.Range("M5").Value = "[5.417132802569313, 53.37944851594434]"
$ws.Range("N5").Value = 0.01740816640052789
$ws.Range("O5").Value = 0.01740816640052789
$ws.Range("P5").Value = 0.5975001042863477
$ws.Range("Q5").Value = "[-0.1823947686768843, 1.3773949772495797]"
$ws.Range("R5").Value = 0.1298193195262545
$ws.Range("S5").Value = 0.1298193195262545
$ws.Range("T5").Value = 61.70189188404225
$ws.Range("U5").Value = "[48.80029224126932, 74.60349152681519]"
$ws.Range("V5").Value = 0.000000000001658451154185059
$ws.Range("W5").Value = 0.000000000001658451154185059
$ws.Range("X5").Value = 22.74930930930975
$ws.Range("Y5").Value = 19.62882882882921
$ws.Range("Z5").Value = 25.86978978979029

# Row 6
$ws.Range("F6").Value = 24.23000000000035
$ws.Range("H6").Value = 0.000008172932638661479
$ws.Range("I6").Value = 0.000008172932638661479
$ws.Range("L6").Value = 48.82383785672109
$ws.Range("M6").Value = "[27.61112473011437, 70.03655098332781]"
$ws.Range("N6").Value = 0.00003060780112718042
$ws.Range("O6").Value = 0.00003060780112718042
$ws.Range("P6").Value = 0.2452895164965003
$ws.Range("Q6").Value = "[-0.27044741562434726, 0.7610264486173479]"
$ws.Range("R6").Value = 0.3432179987819299
$ws.Range("S6").Value = 0.3432179987819299
$ws.Range("T6").Value = 62.46337300229482
$ws.Range("U6").Value = "[49.983144129008025, 74.9436018755816]"
$ws.Range("V6").Value = 0.000000000000404343225568482
$ws.Range("W6").Value = 0.000000000000404343225568482
$ws.Range("X6").Value = 23.28408408408442
$ws.Range("Y6").Value = 21.29523523523554
$ws.Range("Z6").Value = 25.2729329329333

# Row 7
$ws.Range("F7").Value = 24.23000000000035
$ws.Range("H7").Value = 0.0006384567581843026
$ws.Range("I7").Value = 0.0006384567581843026
$ws.Range("L7").Value = 35.754493620107
$ws.Range("M7").Value = "[12.913240167744313, 58.595747072469685]"
$ws.Range("N7").Value = 0.002877647959147689
$ws.Range("O7").Value = 0.002877647959147689
$ws.Range("P7").Value = 0.5346053564667317
$ws.Range("Q7").Value = "[-0.10692107129334616, 1.1761317842268095]"
$ws.Range("R7").Value = 0.1001977283572195
$ws.Range("S7").Value = 0.1001977283572195
$ws.Range("T7").Value = 52.26939996954284
$ws.Range("U7").Value = "[40.118101708064515, 64.42069823102116]"
$ws.Range("V7").Value = 0.00000000003804423442943516
$ws.Range("W7").Value = 0.00000000003804423442943516
$ws.Range("X7").Value = 22.16838838838871
$ws.Range("Y7").Value = 19.69445445445474
$ws.Range("Z7").Value = 24.64232232232268

# Row 8
$ws.Range("F8").Value = 24.23000000000035
$ws.Range("H8").Value = 0.00005114066639166914
$ws.Range("I8").Value = 0.00005114066639166914
$ws.Range("L8").Value = 50.89442197785667
$ws.Range("M8").Value = "[27.39100843474621, 74.39783552096712]"
$ws.Range("N8").Value = 0.0000744269601442582
$ws.Range("O8").Value = 0.0000744269601442582
$ws.Range("P8").Value = 0.1320789704211931
$ws.Range("Q8").Value = "[-0.4339737599553457, 0.6981317007977319]"
$ws.Range("R8").Value = 0.6406550110187896
$ws.Range("S8").Value = 0.6406550110187896
$ws.Range("T8").Value = 66.15439897350153
$ws.Range("U8").Value = "[51.69163160565604, 80.61716634134702]"
$ws.Range("V8").Value = 0.000000000006360689752682447
$ws.Range("W8").Value = 0.000000000006360689752682447
$ws.Range("X8").Value = 23.720660660661
$ws.Range("Y8").Value = 21.53777777777809
$ws.Range("Z8").Value = 25.90354354354391

